{"js": "// Update the benchmark summary table:\n//  - three \"unit\" rows (throughput totals) become \"0M\" placeholders\n//  - their original numeric values get relocated to the bottom three\n//    per-run detail rows, which also collapse from a tab-separated\n//    9-column breakdown down to that single total figure\n//  - a handful of other single-value cells are refreshed with new\n//    measurements\n\nconst body = context.document.body;\n\n// --- 1) Simple single-run numeric replacements (each old value is\n//        unique in the document, so search+replace keeps the existing\n//        run's formatting, e.g. Times New Roman / sz 22, intact). ---\nconst simpleReplacements = [\n  [\"64.62\", \"0M\"],\n  [\"386.08\", \"0M\"],\n  [\"1091\", \"0M\"],\n  [\"5991\", \"6669\"],\n  [\"0.05684\", \"0.06777\"],\n  [\"0.07009\", \"0.07769\"],\n  [\"323.71821\", \"386.07718\"],\n];\n\nfor (const [oldVal, newVal] of simpleReplacements) {\n  const results = body.search(oldVal, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  results.items[0].insertText(newVal, \"Replace\");\n  await context.sync();\n}\n\n// --- 2) Collapse the three trailing detail rows (tab-separated\n//        9-value breakdowns) down to a single value, reusing the\n//        original totals that used to live in rows 0-2. ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\nconst table = tables.items[0];\ntable.load(\"rowCount\");\nawait context.sync();\n\nconst rowCount = table.rowCount;\nconst collapsedRowValues = [\n  [\"64.62\"],\n  [\"386.08\"],\n  [\"1091\"],\n];\n\nfor (let i = 0; i < collapsedRowValues.length; i++) {\n  const rowIndex = rowCount - collapsedRowValues.length + i;\n  const cell = table.getCell(rowIndex, 0);\n  const paragraph = cell.body.paragraphs.getFirst();\n  const range = paragraph.getRange();\n  range.insertText(collapsedRowValues[i][0], \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Update the benchmark summary table:\n#  - three \"unit\" rows (throughput totals) become \"0M\" placeholders\n#  - their original numeric values get relocated to the bottom three\n#    per-run detail rows, which also collapse from a tab-separated\n#    9-column breakdown down to that single total figure\n#  - a handful of other single-value cells are refreshed with new\n#    measurements\n\n$d = $word.ActiveDocument\n\n# --- 1) Simple single-run numeric replacements (each old value is\n#        unique in the document, so Find/Replace keeps the existing\n#        run's formatting, e.g. Times New Roman / sz 22, intact). ---\n$simpleReplacements = @(\n    @{ Old = \"64.62\";      New = \"0M\" },\n    @{ Old = \"386.08\";     New = \"0M\" },\n    @{ Old = \"1091\";       New = \"0M\" },\n    @{ Old = \"5991\";       New = \"6669\" },\n    @{ Old = \"0.05684\";    New = \"0.06777\" },\n    @{ Old = \"0.07009\";    New = \"0.07769\" },\n    @{ Old = \"323.71821\";  New = \"386.07718\" }\n)\n\nforeach ($rep in $simpleReplacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $rep.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $rep.New\n    $find.Execute($rep.Old, $true, $false, $false, $false, $false, $true, 1, $false, $rep.New, 2) | Out-Null\n}\n\n# --- 2) Collapse the three trailing detail rows (tab-separated\n#        9-value breakdowns) down to a single value, reusing the\n#        original totals that used to live in rows 1-3 (1-indexed). ---\n$t = $d.Tables(1)\n$rowCount = $t.Rows.Count\n\n$collapsedRowValues = @(\"64.62\", \"386.08\", \"1091\")\nfor ($i = 0; $i -lt $collapsedRowValues.Count; $i++) {\n    $rowIndex = $rowCount - $collapsedRowValues.Count + 1 + $i\n    $cell = $t.Cell($rowIndex, 1)\n    $r = $cell.Range\n    $r.End = $r.End - 1\n    $r.Text = $collapsedRowValues[$i]\n}\n"}
